# Update cryptocurrency price/volume figures per latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.194.07'
$ws.Range('E2').Value = '  +1.69%  '
$ws.Range('D3').Value = '1.785.13'
$ws.Range('E3').Value = '  +1.01%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.90'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.29%  '
$ws.Range('E6').Value = '  +0.66%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.83'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.46%  '
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0687'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.35%  '
$ws.Range('E11').Value = '  +1.11%  '
$ws.Range('D12').Value = '2.043.07'
$ws.Range('E12').Value = '  +1.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.04'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.40%  '
$ws.Range('D14').Value = '1.789.70'
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '34.152.66'
$ws.Range('E15').Value = '  +1.44%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.623'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.82%  '
$ws.Range('E17').Value = '  +1.80%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.09'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '246.29'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.92%  '
$ws.Range('E20').Value = '  +0.76%  '
$ws.Range('E21').Value = '  +4.19%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('E23').Value = '  +2.67%  '
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '161.45'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('E26').Value = '  +2.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.31'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.65%  '
$ws.Range('E28').Value = '  +2.10%  '
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.23'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.43%  '
$ws.Range('E31').Value = '  +1.82%  '
$ws.Range('E32').Value = '  +3.13%  '
$ws.Range('E33').Value = '  +4.44%  '
$ws.Range('E34').Value = '  +1.04%  '
$ws.Range('D35').Value = '1.446.77'
$ws.Range('E35').Value = '  +5.10%  '
$ws.Range('E36').Value = '  +2.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.44'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +10.20%  '
$ws.Range('E38').Value = '  +4.47%  '
$ws.Range('E39').Value = '  +1.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '80.25'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.37'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.926'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.68'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.49'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.75%  '
$ws.Range('E45').Value = '  +4.75%  '
$ws.Range('E46').Value = '  +2.04%  '
$ws.Range('E47').Value = '  -0.54%  '
$ws.Range('D48').Value = '0.0₆0135'
$ws.Range('E48').Value = '  -1.11%  '
$ws.Range('D49').Value = '1.944.97'
$ws.Range('E49').Value = '  +1.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '106.02'
$ws.Range('D50').ClearFormats()
$ws.Range('E51').Value = '  +0.06%  '
